$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Will" -> "Will/Laurance" (existing merged cell G11:H11, style already applied) ---
$ws.Range("G11").Value = "Will/Laurance"

# --- Add new rows 19-21 in column B (pushing the former row 21 content down to make room) ---
# Row 21 (was row 21 before the insert, now additional content alongside the existing D21/E21 block)
$ws.Range("B21").Value = "Transfer learning"

# Row 19: bold label
$ws.Range("B19").Value = "To look into"
$ws.Range("B19").Font.Bold = $true

# Row 20: wrapped text, taller row
$ws.Range("B20").Value = "filtering and use of Classical methods"
$ws.Range("B20").WrapText = $true
$ws.Range("B20").RowHeight = 45

# --- Selection / zoom bookkeeping to mirror the author's view state ---
$excel.ActiveWindow.Zoom = 205
[void]$ws.Range("E16").Select()
